# Updated cryptos list values (Price + Volume(1h)) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.382.17"
Set-TextValue $ws.Range("E2") "  +1.42%  "
Set-TextValue $ws.Range("D3") "1.828.46"
Set-TextValue $ws.Range("E3") "  +0.34%  "
Set-TextValue $ws.Range("D4") "0.9993"
Set-TextValue $ws.Range("E4") "  -0.27%  "
Set-TextValue $ws.Range("D5") "313.81"
Set-TextValue $ws.Range("E5") "  +0.78%  "
Set-TextValue $ws.Range("D6") "0.9992"
Set-TextValue $ws.Range("E6") "  -0.24%  "
Set-TextValue $ws.Range("D7") "0.4491"
Set-TextValue $ws.Range("E7") "  +5.91%  "
Set-TextValue $ws.Range("D8") "0.3781"
Set-TextValue $ws.Range("E8") "  +3.30%  "
Set-TextValue $ws.Range("D9") "0.07497"
Set-TextValue $ws.Range("E9") "  +3.92%  "
Set-TextValue $ws.Range("D10") "0.8963"
Set-TextValue $ws.Range("E10") "  +6.67%  "
Set-TextValue $ws.Range("D11") "21.11"
Set-TextValue $ws.Range("E11") "  +2.46%  "
Set-TextValue $ws.Range("D12") "1.831.75"
Set-TextValue $ws.Range("E12") "  +0.42%  "
Set-TextValue $ws.Range("D13") "6.781"
Set-TextValue $ws.Range("E13") "  +2.08%  "
Set-TextValue $ws.Range("D14") "94.70"
Set-TextValue $ws.Range("E14") "  +5.71%  "
Set-TextValue $ws.Range("D15") "5.420"
Set-TextValue $ws.Range("E15") "  +2.75%  "
Set-TextValue $ws.Range("D16") "0.07125"
Set-TextValue $ws.Range("E16") "  +1.20%  "
Set-TextValue $ws.Range("D17") "1.0000"
Set-TextValue $ws.Range("E17") "  -0.31%  "
Set-TextValue $ws.Range("D18") "0.000008833"
Set-TextValue $ws.Range("E18") "  +1.27%  "
Set-TextValue $ws.Range("D19") "0.9996"
Set-TextValue $ws.Range("E19") "  -0.16%  "
Set-TextValue $ws.Range("D20") "15.24"
Set-TextValue $ws.Range("E20") "  +2.69%  "
Set-TextValue $ws.Range("D21") "27.381.94"
Set-TextValue $ws.Range("E21") "  +0.95%  "
Set-TextValue $ws.Range("D22") "5.299"
Set-TextValue $ws.Range("E22") "  +3.46%  "
Set-TextValue $ws.Range("D23") "10.99"
Set-TextValue $ws.Range("E23") "  +1.94%  "
Set-TextValue $ws.Range("D24") "2.051.38"
Set-TextValue $ws.Range("E24") "  -0.18%  "
Set-TextValue $ws.Range("D25") "2.004"
Set-TextValue $ws.Range("E25") "  +1.28%  "
Set-TextValue $ws.Range("D26") "2.487"
Set-TextValue $ws.Range("E26") "  +10.92%  "
Set-TextValue $ws.Range("D27") "152.00"
Set-TextValue $ws.Range("E27") "  +0.40%  "
Set-TextValue $ws.Range("D28") "18.64"
Set-TextValue $ws.Range("E28") "  +2.68%  "
Set-TextValue $ws.Range("D29") "5.397"
Set-TextValue $ws.Range("E29") "  +3.45%  "
Set-TextValue $ws.Range("D30") "118.27"
Set-TextValue $ws.Range("E30") "  +1.48%  "
Set-TextValue $ws.Range("D31") "0.08848"
Set-TextValue $ws.Range("E31") "  +1.45%  "
Set-TextValue $ws.Range("D32") "0.7825"
Set-TextValue $ws.Range("E32") "  +6.82%  "
Set-TextValue $ws.Range("D33") "1.201"
Set-TextValue $ws.Range("E33") "  +2.37%  "
Set-TextValue $ws.Range("D34") "4.594"
Set-TextValue $ws.Range("E34") "  +4.27%  "
Set-TextValue $ws.Range("D35") "2.885"
Set-TextValue $ws.Range("E35") "  -0.63%  "
Set-TextValue $ws.Range("D36") "0.9986"
Set-TextValue $ws.Range("E36") "  -0.28%  "
Set-TextValue $ws.Range("D37") "1.112"
Set-TextValue $ws.Range("E37") "  +2.00%  "
Set-TextValue $ws.Range("D38") "0.01993"
Set-TextValue $ws.Range("E38") "  +2.66%  "
Set-TextValue $ws.Range("D39") "0.05343"
Set-TextValue $ws.Range("E39") "  +2.52%  "
Set-TextValue $ws.Range("D40") "7.428"
Set-TextValue $ws.Range("E40") "  +2.41%  "
Set-TextValue $ws.Range("D41") "0.5361"
Set-TextValue $ws.Range("E41") "  +5.19%  "
Set-TextValue $ws.Range("D42") "0.1735"
Set-TextValue $ws.Range("E42") "  +3.04%  "
Set-TextValue $ws.Range("D43") "2.858"
Set-TextValue $ws.Range("E43") "  -0.37%  "
Set-TextValue $ws.Range("D44") "2.280"
Set-TextValue $ws.Range("E44") "  +17.18%  "
Set-TextValue $ws.Range("D45") "8.849"
Set-TextValue $ws.Range("E45") "  +3.83%  "
Set-TextValue $ws.Range("D46") "0.5160"
Set-TextValue $ws.Range("E46") "  +9.35%  "
Set-TextValue $ws.Range("D47") "10.77"
Set-TextValue $ws.Range("E47") "  +2.80%  "
Set-TextValue $ws.Range("D48") "106.80"
Set-TextValue $ws.Range("E48") "  +1.40%  "
Set-TextValue $ws.Range("D49") "1.709"
Set-TextValue $ws.Range("E49") "  +3.98%  "
Set-TextValue $ws.Range("D50") "0.9984"
Set-TextValue $ws.Range("E50") "  -0.24%  "
Set-TextValue $ws.Range("D51") "0.06382"
Set-TextValue $ws.Range("E51") "  +0.97%  "
